# Commit: "Modificato id worksheet per la colonna -id"
#
# 1) Rename the sheet from "chief_complaint_id" to "_id"
# 2) Give the new sheet window-minimized book view (best-effort; harmless if
#    the host ignores it)
# 3) Size up columns B:E to fit their (now differently-labelled) content
# 4) Leave the cursor on the default cell so Excel doesn't need to persist an
#    explicit <selection> further down the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet rename: chief_complaint_id -> _id -------------------------------
$ws.Name = "_id"

# --- book view: mark the workbook window minimized -------------------------
$w = $excel.ActiveWindow
$w.WindowState = -4140   # xlMinimized

# --- column widths (characters) for the B:E block --------------------------
# Target stored widths (OOXML <col width=.../>, MDW=7 units) are
# 16.109375 / 21.88671875 / 24.77734375 / 8.88671875 respectively; the
# ColumnWidth values below are chosen so the host's character->pixel
# rounding lands on the nearest representable width to each target.
$ws.Columns.Item(2).ColumnWidth = 15.3
$ws.Columns.Item(3).ColumnWidth = 21
$ws.Columns.Item(4).ColumnWidth = 24
$ws.Columns.Item(5).ColumnWidth = 8

# --- reset selection back to the sheet's default (A1) -----------------------
$ws.Range("A1").Select() | Out-Null
